$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header labels in C1 and D1: Mean/Median Absolute Error were reversed.
$ws.Range("C1").Value = "Mean Absolute Error"
$ws.Range("D1").Value = "Median Absolute Error"
